$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '30.788.69'
$ws.Range("E2").Value = '  +1.73%  '

$ws.Range("D3").Value = '2.114.21'
$ws.Range("E3").Value = '  +6.52%  '

$ws.Range("E4").Value = '  +0.30%  '

$ws.Range("D5").Value = "'333.20"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +3.21%  '

$ws.Range("E6").Value = '  +0.31%  '

$ws.Range("D7").Value = "'0.5326"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +3.97%  '

$ws.Range("D8").Value = "'0.4387"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +6.93%  '

$ws.Range("D9").Value = "'0.09021"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +6.78%  '

$ws.Range("D10").Value = "'46.03"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +8.15%  '

$ws.Range("D11").Value = "'1.180"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +4.06%  '

$ws.Range("D12").Value = "'25.02"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +2.95%  '

$ws.Range("D13").Value = '2.107.57'
$ws.Range("E13").Value = '  +7.37%  '

$ws.Range("D14").Value = "'6.753"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +4.48%  '

$ws.Range("D15").Value = "'7.816"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +5.62%  '

$ws.Range("D16").Value = "'97.43"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +3.85%  '

$ws.Range("E17").Value = '  +0.51%  '

$ws.Range("E18").Value = '  +2.06%  '

$ws.Range("D19").Value = "'0.06659"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +1.89%  '

$ws.Range("D20").Value = "'19.13"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +2.05%  '

$ws.Range("E21").Value = '  +0.32%  '

$ws.Range("D22").Value = "'6.359"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +4.84%  '

$ws.Range("D23").Value = '30.844.73'
$ws.Range("E23").Value = '  +1.72%  '

$ws.Range("D24").Value = "'12.36"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +7.82%  '

$ws.Range("D25").Value = '2.355.77'
$ws.Range("E25").Value = '  +7.51%  '

$ws.Range("D26").Value = "'2.252"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +2.43%  '

$ws.Range("E27").Value = '  +0.61%  '

$ws.Range("E28").Value = '  +8.27%  '

$ws.Range("D29").Value = "'163.06"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +0.24%  '

$ws.Range("D30").Value = "'133.17"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +2.06%  '

$ws.Range("D31").Value = "'1.170"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +2.98%  '

$ws.Range("E32").Value = '  +2.08%  '

$ws.Range("D33").Value = "'6.225"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +3.21%  '

$ws.Range("E34").Value = '  +5.81%  '

$ws.Range("D35").Value = "'1.544"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +17.43%  '

$ws.Range("D36").Value = "'0.02611"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +5.45%  '

$ws.Range("D37").Value = "'5.534"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +3.06%  '

$ws.Range("D38").Value = "'9.558"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +7.33%  '

$ws.Range("B39").Value = 'Aptos'
$ws.Range("C39").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D39").Value = "'12.85"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +9.47%  '

$ws.Range("B40").Value = 'Hedera'
$ws.Range("C40").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D40").Value = "'0.06744"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +3.77%  '

$ws.Range("D41").Value = "'0.2283"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +5.03%  '

$ws.Range("D42").Value = "'0.6870"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +4.34%  '

$ws.Range("D43").Value = "'1.252"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +2.42%  '

$ws.Range("B44").Value = 'Decentraland'
$ws.Range("C44").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range("D44").Value = "'0.6454"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +5.46%  '

$ws.Range("B45").Value = 'Frax'
$ws.Range("C45").Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range("D45").Value = "'1.001"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +0.63%  '

$ws.Range("B46").Value = 'EnergySwap'
$ws.Range("C46").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D46").Value = "'14.11"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +4.30%  '

$ws.Range("D47").Value = "'2.232"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +2.00%  '

$ws.Range("D48").Value = "'3.672"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +0.86%  '

$ws.Range("D49").Value = "'1.277"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +4.38%  '

$ws.Range("D50").Value = "'82.53"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +3.77%  '

$ws.Range("D51").Value = "'120.84"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -2.20%  '
